$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "DESIRED PUBLICATION DATE:  August 3, 2017" -> "... August 2, 2017"
# ------------------------------------------------------------------
$d.Content.Find.Execute("August 3, ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "August 2, ", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "...in the front meeting room of the Crosspoint United Methodist
#     Church..." -> "...in Room 400 of the Crosspoint United Methodist
#     Church..."
#    The _GoBack bookmark (originally right after "Church") needs to
#    end up right after the new "Room 400" text instead.
# ------------------------------------------------------------------
$d.Content.Find.Execute("the front meeting room", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Room 400", 2) | Out-Null

$roomRange = $d.Content
$roomRange.Find.Execute("Room 400", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null

$d.Bookmarks.Add("_GoBack", $d.Range($roomRange.End, $roomRange.End)) | Out-Null

$d.Save()
